# Weekly data refresh: insert a new data row for "Ciboulette" (Hortaliza)
# at the top of the existing data block (row 529), pushing the remaining
# rows (529-598) down by one (to 530-599).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 529, shifting rows 529:598 down to 530:599.
$ws.Rows.Item(529).Insert()

# Populate the newly inserted row with the new week's record.
$ws.Range("A529").Value = 9
$ws.Range("B529").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C529").Value = 'Metropolitana'
$ws.Range("D529").Value2 = 45127
$ws.Range("E529").Value = 13
$ws.Range("F529").Value = 100112039
$ws.Range("G529").Value = 'Ciboulette'
$ws.Range("H529").Value = 'Sin especificar'
$ws.Range("I529").Value = 'Primera'
$ws.Range("J529").Value = 340
$ws.Range("K529").Value = 1500
$ws.Range("L529").Value = 1500
$ws.Range("M529").Value = 1500
$ws.Range("N529").Value = '$/docena de atados'
$ws.Range("O529").Value = 'Región Metropolitana'
$ws.Range("P529").Value = 500
$ws.Range("Q529").Value = 3
$ws.Range("R529").Value = 'Hortaliza'
